$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the average_doctor / average_doctor_old header labels (BP1/BQ1)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Updated per-app "_old" (Harvard) classification stats and recomputed
# average_doctor / average_doctor_old columns
$ws.Range("E4").Value = 0.473
$ws.Range("F4").Value = 0.053
$ws.Range("G4").Value = 0.229
$ws.Range("N4").Value = 0.477
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.25
$ws.Range("Q4").Value = 0.052
$ws.Range("R4").Value = 0.035
$ws.Range("S4").Value = 0.188
$ws.Range("W4").Value = 0.375
$ws.Range("X4").Value = 0.105
$ws.Range("Y4").Value = 0.324
$ws.Range("AI4").Value = 0.411
$ws.Range("AJ4").Value = 0.092
$ws.Range("AK4").Value = 0.303
$ws.Range("AU4").Value = 0.245
$ws.Range("AW4").Value = 0.158
$ws.Range("BA4").Value = 2.053
$ws.Range("BB4").Value = 0.147
$ws.Range("BC4").Value = 0.383
$ws.Range("BG4").Value = 0.722
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.75
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.6840000000000001
$ws.Range("BQ4").Value = 0.762
$ws.Range("E5").Value = 0.598
$ws.Range("F5").Value = 0.06
$ws.Range("G5").Value = 0.245
$ws.Range("N5").Value = 0.717
$ws.Range("O5").Value = 0.077
$ws.Range("P5").Value = 0.278
$ws.Range("Q5").Value = 0.035
$ws.Range("R5").Value = 0.014
$ws.Range("S5").Value = 0.12
$ws.Range("W5").Value = 0.343
$ws.Range("X5").Value = 0.097
$ws.Range("Y5").Value = 0.311
$ws.Range("AI5").Value = 0.412
$ws.Range("AJ5").Value = 0.08799999999999999
$ws.Range("AK5").Value = 0.297
$ws.Range("AU5").Value = 0.456
$ws.Range("AV5").Value = 0.074
$ws.Range("AW5").Value = 0.272
$ws.Range("BA5").Value = 1.298
$ws.Range("BB5").Value = 0.074
$ws.Range("BC5").Value = 0.271
$ws.Range("BG5").Value = 0.375
$ws.Range("BH5").Value = 0.049
$ws.Range("BI5").Value = 0.221
$ws.Range("BM5").Value = 0.525
$ws.Range("BN5").Value = 0.047
$ws.Range("BO5").Value = 0.217
$ws.Range("BP5").Value = 0.433
$ws.Range("BQ5").Value = 0.458
$ws.Range("E6").Value = 0.528
$ws.Range("N6").Value = 0.573
$ws.Range("Q6").Value = 0.042
$ws.Range("W6").Value = 0.358
$ws.Range("AI6").Value = 0.411
$ws.Range("AU6").Value = 0.319
$ws.Range("BA6").Value = 1.584
$ws.Range("BG6").Value = 0.494
$ws.Range("BM6").Value = 0.618
$ws.Range("BP6").Value = 0.528
$ws.Range("BQ6").Value = 0.57
$ws.Range("E7").Value = 0.5679999999999999
$ws.Range("N7").Value = 0.651
$ws.Range("Q7").Value = 0.037
$ws.Range("W7").Value = 0.349
$ws.Range("AI7").Value = 0.412
$ws.Range("AU7").Value = 0.389
$ws.Range("BA7").Value = 1.399
$ws.Range("BG7").Value = 0.415
$ws.Range("BM7").Value = 0.5590000000000001
$ws.Range("BP7").Value = 0.466
$ws.Range("BQ7").Value = 0.497
$ws.Range("E8").Value = 0.704
$ws.Range("F8").Value = 0.07199999999999999
$ws.Range("G8").Value = 0.268
$ws.Range("N8").Value = 0.802
$ws.Range("O8").Value = 0.06
$ws.Range("P8").Value = 0.244
$ws.Range("Q8").Value = 0.038
$ws.Range("W8").Value = 0.414
$ws.Range("X8").Value = 0.118
$ws.Range("Y8").Value = 0.343
$ws.Range("AI8").Value = 0.483
$ws.Range("AJ8").Value = 0.135
$ws.Range("AK8").Value = 0.367
$ws.Range("AU8").Value = 0.403
$ws.Range("AV8").Value = 0.081
$ws.Range("AW8").Value = 0.285
$ws.Range("BA8").Value = 1.769
$ws.Range("BB8").Value = 0.109
$ws.Range("BC8").Value = 0.33
$ws.Range("BG8").Value = 0.5590000000000001
$ws.Range("BH8").Value = 0.11
$ws.Range("BI8").Value = 0.331
$ws.Range("BM8").Value = 0.6830000000000001
$ws.Range("BN8").Value = 0.059
$ws.Range("BO8").Value = 0.243
$ws.Range("BP8").Value = 0.59
$ws.Range("BQ8").Value = 0.625
$ws.Range("E9").Value = 0.667
$ws.Range("F9").Value = 0.222
$ws.Range("G9").Value = 0.471
$ws.Range("N9").Value = 0.738
$ws.Range("O9").Value = 0.193
$ws.Range("P9").Value = 0.44
$ws.Range("W9").Value = 0.31
$ws.Range("X9").Value = 0.214
$ws.Range("Y9").Value = 0.462
$ws.Range("AI9").Value = 0.429
$ws.Range("AJ9").Value = 0.245
$ws.Range("AK9").Value = 0.495
$ws.Range("BA9").Value = 1.738
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.595
$ws.Range("BH9").Value = 0.241
$ws.Range("BI9").Value = 0.491
$ws.Range("BM9").Value = 0.667
$ws.Range("BN9").Value = 0.222
$ws.Range("BO9").Value = 0.471
$ws.Range("BP9").Value = 0.579
$ws.Range("BQ9").Value = 0.617
$ws.Range("E10").Value = 0.8100000000000001
$ws.Range("F10").Value = 0.154
$ws.Range("G10").Value = 0.393
$ws.Range("N10").Value = 0.929
$ws.Range("O10").Value = 0.066
$ws.Range("P10").Value = 0.258
$ws.Range("W10").Value = 0.524
$ws.Range("X10").Value = 0.249
$ws.Range("Y10").Value = 0.499
$ws.Range("AI10").Value = 0.524
$ws.Range("AJ10").Value = 0.249
$ws.Range("AK10").Value = 0.499
$ws.Range("AU10").Value = 0.405
$ws.Range("AV10").Value = 0.241
$ws.Range("AW10").Value = 0.491
$ws.Range("BA10").Value = 2.19
$ws.Range("BB10").Value = 0.214
$ws.Range("BC10").Value = 0.462
$ws.Range("BG10").Value = 0.667
$ws.Range("BH10").Value = 0.222
$ws.Range("BI10").Value = 0.471
$ws.Range("BM10").Value = 0.833
$ws.Range("BN10").Value = 0.139
$ws.Range("BO10").Value = 0.373
$ws.Range("BP10").Value = 0.73
$ws.Range("BQ10").Value = 0.759
$ws.Range("E11").Value = 0.857
$ws.Range("F11").Value = 0.122
$ws.Range("G11").Value = 0.35
$ws.Range("N11").Value = 0.929
$ws.Range("O11").Value = 0.066
$ws.Range("P11").Value = 0.258
$ws.Range("W11").Value = 0.524
$ws.Range("X11").Value = 0.249
$ws.Range("Y11").Value = 0.499
$ws.Range("AI11").Value = 0.595
$ws.Range("AJ11").Value = 0.241
$ws.Range("AK11").Value = 0.491
$ws.Range("AU11").Value = 0.571
$ws.Range("AV11").Value = 0.245
$ws.Range("AW11").Value = 0.495
$ws.Range("BA11").Value = 2.19
$ws.Range("BB11").Value = 0.214
$ws.Range("BC11").Value = 0.462
$ws.Range("BG11").Value = 0.667
$ws.Range("BH11").Value = 0.222
$ws.Range("BI11").Value = 0.471
$ws.Range("BM11").Value = 0.833
$ws.Range("BN11").Value = 0.139
$ws.Range("BO11").Value = 0.373
$ws.Range("BP11").Value = 0.73
$ws.Range("BQ11").Value = 0.766
$ws.Range("E12").Value = 1.417
$ws.Range("F12").Value = 0.854
$ws.Range("G12").Value = 0.924
$ws.Range("N12").Value = 1.256
$ws.Range("O12").Value = 0.293
$ws.Range("P12").Value = 0.542
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.432
$ws.Range("Y12").Value = 0.657
$ws.Range("AI12").Value = 1.6
$ws.Range("AJ12").Value = 1.44
$ws.Range("AK12").Value = 1.2
$ws.Range("AU12").Value = 2.846
$ws.Range("AV12").Value = 3.361
$ws.Range("AW12").Value = 1.833
$ws.Range("BA12").Value = 3.786
$ws.Range("BB12").Value = 0.449
$ws.Range("BC12").Value = 0.67
$ws.Range("BG12").Value = 1.143
$ws.Range("BH12").Value = 0.194
$ws.Range("BI12").Value = 0.44
$ws.Range("BM12").Value = 1.229
$ws.Range("BN12").Value = 0.233
$ws.Range("BO12").Value = 0.483
$ws.Range("BP12").Value = 1.262
$ws.Range("BQ12").Value = 1.246
$ws.Range("E13").Value = 1.415
$ws.Range("F13").Value = 0.295
$ws.Range("G13").Value = 0.543
$ws.Range("N13").Value = 1.73
$ws.Range("O13").Value = 0.476
$ws.Range("P13").Value = 0.6899999999999999
$ws.Range("W13").Value = 0.985
$ws.Range("X13").Value = 0.194
$ws.Range("Y13").Value = 0.441
$ws.Range("AI13").Value = 1.154
$ws.Range("AJ13").Value = 0.303
$ws.Range("AK13").Value = 0.551
$ws.Range("AU13").Value = 2.039
$ws.Range("AV13").Value = 0.339
$ws.Range("AW13").Value = 0.582
$ws.Range("BA13").Value = 2.171
$ws.Range("BB13").Value = 0.283
$ws.Range("BC13").Value = 0.532
$ws.Range("BG13").Value = 0.539
$ws.Range("BH13").Value = 0.051
$ws.Range("BI13").Value = 0.226
$ws.Range("BM13").Value = 0.787
$ws.Range("BN13").Value = 0.163
$ws.Range("BO13").Value = 0.403
$ws.Range("BP13").Value = 0.724
$ws.Range("BQ13").Value = 0.667
